# Applies the header-rename + column-width tweaks described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes -------------------------------------------------
# NOTE: the COM layer stores column width internally in pixel-rounded units
# and re-derives the OOXML "characters" width with a +0.8333 fudge factor,
# so asking for the exact integer width we want (e.g. 31) actually persists
# as 31.8333 in the saved file. Biasing the requested width down by 0.875
# lands safely inside the pixel bucket that round-trips to the desired
# whole-number width in the saved <col width="..."/> attribute.
$bias = 0.875
$ws.Columns.Item(2).ColumnWidth = 31 - $bias    # B: 37 -> 31
$ws.Columns.Item(15).ColumnWidth = 46 - $bias   # O: 44 -> 46
$ws.Columns.Item(16).ColumnWidth = 48 - $bias   # P: 46 -> 48
$ws.Columns.Item(17).ColumnWidth = 48 - $bias   # Q: 46 -> 48
$ws.Columns.Item(42).ColumnWidth = 27 - $bias   # AP: 29 -> 27
$ws.Columns.Item(43).ColumnWidth = 29 - $bias   # AQ: 31 -> 29
$ws.Columns.Item(44).ColumnWidth = 30 - $bias   # AR: 32 -> 30
$ws.Columns.Item(45).ColumnWidth = 32 - $bias   # AS: 34 -> 32
$ws.Columns.Item(46).ColumnWidth = 27 - $bias   # AT: 29 -> 27
$ws.Columns.Item(47).ColumnWidth = 29 - $bias   # AU: 31 -> 29
$ws.Columns.Item(48).ColumnWidth = 33 - $bias   # AV: 30 -> 33
$ws.Columns.Item(49).ColumnWidth = 26 - $bias   # AW: 23 -> 26

# --- Header text renames (row 1) ------------------------------------------
$ws.Range("A1").Value = "button_alertActions_class"
$ws.Range("B1").Value = "div_testCaseData_internalText"
$ws.Range("N1").Value = "div_thumbnailImages_class"
$ws.Range("O1").Value = "header_testRunTitles_internalRoleHeadingName"
$ws.Range("P1").Value = "header_testRunTitles_internalRoleHeadingName_1"
$ws.Range("Q1").Value = "header_testRunTitles_internalRoleHeadingName_2"
$ws.Range("AP1").Value = "link_testRunLinks_plan_id"
$ws.Range("AQ1").Value = "link_testRunLinks_plan_id_1"
$ws.Range("AR1").Value = "link_testRunLinks_project_id"
$ws.Range("AS1").Value = "link_testRunLinks_project_id_1"
$ws.Range("AT1").Value = "link_testRunLinks_team_id"
$ws.Range("AU1").Value = "link_testRunLinks_team_id_1"
$ws.Range("AV1").Value = "span_logDetails_internalHasText"
$ws.Range("AW1").Value = "span_logDetails_nthChild"
